$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header for column F ---
$ws.Range("F1").Value = "Farmplot price"

# --- New farmplot-price formulas in column F (rows 2-10) ---
$ws.Range("F2").Formula = "=ROUNDDOWN(`$K`$3*`$K`$4^(A2-1),0)"
$ws.Range("F3:F10").Formula = "=ROUNDDOWN(`$K`$3*`$K`$4^(A3-1),0)"

# --- New parameter labels/values used by the formula above ---
$ws.Range("J3").Value = "Base farmplot price"
$ws.Range("K3").Value = 5000

$ws.Range("J4").Value = "Plot price multi"
$ws.Range("K4").Value = 1.5

# --- Column widths (closest representable values on this host's width grid) ---
$ws.Columns.Item(6).ColumnWidth = 13.333333333333334
$ws.Columns.Item(10).ColumnWidth = 17.833333333333332

# --- Selection moves from D2 to F2 ---
$ws.Range("F2").Select()
